# Bump the auto-updating "Date" footer placeholder shown on the
# Slide Master and every Slide Layout from 14/02/2025 to 15/02/2025.
#
# ppPlaceholderDate == 16 in the PpPlaceholderType enum, so we find the
# shape that way (robust against the placeholder's differing index/name
# across the various layouts) rather than hard-coding shape indices.

$p = $ppt.ActivePresentation

$oldDateText = "14/02/2025"
$newDateText = "15/02/2025"
$ppPlaceholderDate = 16

function Update-DateShape {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)

        $placeholderType = $null
        if ($shape.Type -eq 14) {
            # msoPlaceholder
            $placeholderType = $shape.PlaceholderFormat.Type
        }

        if ($placeholderType -eq $ppPlaceholderDate -and $shape.HasTextFrame) {
            $textRange = $shape.TextFrame.TextRange
            if ($textRange.Text -eq $oldDateText) {
                $textRange.Text = $newDateText
            }
        }
    }
}

# The Date placeholder on the Slide Master itself.
Update-DateShape $p.SlideMaster.Shapes

# The Date placeholder repeated on every Slide Layout under the master.
$master = $p.SlideMaster
for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Update-DateShape $layout.Shapes
}
